$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $found = $d.Content.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $found) {
        Write-Host "NOT FOUND: $old"
    }
}

# --- Title ---
Replace-Text "The Enigmatic Cosmos: Unveiling Dark Energy's Secrets" "The Profound Evolution of Artistic Expression: Embracing Diversity and Contemporary Forms"

# --- Author name: "Dr. Elara Vanderstelt" (3 runs) -> "Isabella Kingsley" (1 run) ---
Replace-Text "Dr. Elara Vanderstelt" "Isabella Kingsley"

# --- Email: "evanderstelt@cosmicmysteries" + "." + "org" -> "isabella" + "." + "kingsley@educationalhaven" + "." + "org" ---
Replace-Text "evanderstelt@cosmicmysteries" "isabella.kingsley@educationalhaven"

# --- Body paragraph 1, sentence group 1 (before first blank line) ---
Replace-Text "Across the vast expanse of the cosmos, amidst shimmering constellations and radiant nebulae, lies a profound enigma that captivates and perplexes scientists: dark energy" "In the ever-evolving realm of artistic expression, the landscape of creativity has been dramatically transformed over time"

Replace-Text " This elusive force permeates the universe, driving its accelerated expansion at an ever-increasing rate" " From the dawn of humanity, art has served as a profound and multifaceted means of communication, allowing individuals to connect with their innermost thoughts, experiences, and surroundings"

Replace-Text " Unraveling the secrets of dark energy holds the key to understanding the ultimate fate of our universe, yet its true nature remains shrouded in mystery. Embark on a cosmic journey as we delve into the depths of this enigmatic phenomenon, exploring its profound implications for our comprehension of space, time, and the fundamental forces that govern our existence" " In this essay, we will explore the captivating journey of artistic evolution, delving into its diverse manifestations and embracing the contemporary forms that continue to push the boundaries of creativity"

# --- Body paragraph 1, sentence group 2 ---
Replace-Text "The existence of dark energy was first hinted at in the late 1990s when astronomers observed that the expansion of the universe was not decelerating as expected, but rather accelerating" "Throughout history, art has undergone remarkable shifts, reflecting the cultural, social, and technological changes that have shaped human civilization"

Replace-Text " This unexpected discovery challenged prevailing theories and ignited a scientific quest to identify the mysterious force responsible" " The Renaissance witnessed a resurgence of classical ideals, while the Baroque period exuded grandeur and flamboyance"

Replace-Text " Enter dark energy, a hypothetical form of energy permeating the entire universe and causing its expansion to accelerate" " Modernism shattered traditional norms, giving rise to bold abstraction and unconventional artistic expressions"

Replace-Text " But what is this enigmatic entity? Is it a cosmological constant, a manifestation of vacuum energy, or something else entirely?" " These eras of artistic evolution have left an indelible mark on our understanding of beauty, creativity, and the human condition."

# --- Body paragraph 1, sentence group 3 ---
Replace-Text "As scientists grapple with these perplexing questions, the study of dark energy has illuminated other cosmic mysteries" "As we venture into the contemporary art scene, we encounter a kaleidoscope of diverse practices that challenge conventional notions of artistic expression"

Replace-Text " For instance, it has provided insights into the geometry of the universe, suggesting it may be flat or even slightly curved" " Street art adorns urban landscapes, conveying powerful messages of social commentary and political resistance"

Replace-Text " Moreover, dark energy's influence on the cosmic microwave background radiation has helped refine our understanding of the universe's composition and evolution" " Digital art harnesses the boundless possibilities of technology, blurring the lines between reality and the virtual realm"

Replace-Text " While these discoveries have expanded our knowledge, they have also deepened the enigma surrounding dark energy, beckoning us to unravel its true identity and uncover its profound implications for the cosmos" " Performance art pushes the boundaries of artistic engagement, inviting audiences to become active participants in the creative process. These contemporary forms reflect the dynamic and interconnected world we inhabit, showcasing the limitless potential of artistic innovation"

# --- Summary paragraph ---
Replace-Text "Dark energy, a mysterious and enigmatic force, permeates the universe, driving its accelerated expansion" "The journey of artistic evolution has been a captivating tapestry of creativity, reflecting the cultural, social, and technological transformations that have shaped human civilization"

Replace-Text " Its existence and nature remain shrouded in uncertainty, challenging prevailing theories and captivating scientists worldwide" " From the Renaissance to the Baroque period, from Modernism to the diverse contemporary forms, art has served as a profound means of communication, allowing individuals to connect with their innermost thoughts, experiences, and surroundings"

Replace-Text " The quest to understand dark energy has illuminated other cosmic mysteries, providing insights into the geometry of the universe and the composition of the cosmic microwave background radiation. However, these discoveries have also intensified the enigma, urging us to unveil the true identity of dark energy and fathom its profound implications, not only for our understanding of the cosmos but also for the ultimate fate of our universe" " As we continue to embrace the ever-evolving landscape of artistic expression, we can appreciate the boundless potential of human imagination and creativity"

# --- Add trailing empty paragraph at the end of the document ---
$last = $d.Paragraphs.Last
$last.Range.InsertParagraphAfter()

Write-Host "All replacements applied."
